$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.465.85"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").Value = "2.520.14"
$ws.Range("E3").Value = "  +2.54%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.46"
$ws.Range("E5").Value = "  +1.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.03"
$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +0.95%  "

$ws.Range("D9").Value = "2.517.64"
$ws.Range("E9").Value = "  +2.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0976"
$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("E11").Value = "  -1.23%  "

$ws.Range("E12").Value = "  -2.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.333"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").Value = "2.964.01"
$ws.Range("E14").Value = "  +2.42%  "

$ws.Range("D15").Value = "58.444.31"
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.14"
$ws.Range("E16").Value = "  +1.18%  "

$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("D18").Value = "2.521.94"
$ws.Range("E18").Value = "  +2.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.66"
$ws.Range("E19").Value = "  +1.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.69"
$ws.Range("E20").Value = "  +1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.16"
$ws.Range("E21").Value = "  +1.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  +8.62%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.76"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.406"
$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  +0.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  +1.72%  "

$ws.Range("E29").Value = "  +2.42%  "

$ws.Range("E30").Value = "  +2.69%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.20"
$ws.Range("E31").Value = "  +3.65%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.79"
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.35"
$ws.Range("E33").Value = "  +2.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("E36").Value = "  +0.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.26"
$ws.Range("E37").Value = "  -4.61%  "

$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("E39").Value = "  +1.94%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.772"
$ws.Range("E41").Value = "  -1.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "278.79"
$ws.Range("E42").Value = "  +3.16%  "

$ws.Range("E43").Value = "  +2.72%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.00"
$ws.Range("E44").Value = "  +1.53%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "130.08"
$ws.Range("E45").Value = "  +5.55%  "

$ws.Range("E46").Value = "  +1.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0922"
$ws.Range("E47").Value = "  +1.94%  "

$ws.Range("E48").Value = "  +3.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.74"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("E50").Value = "  +1.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.89"
$ws.Range("E51").Value = "  +1.27%  "
